# "update scripts wuth new tpm"
#
# NATMI LR-pairs output for Wnt2 -> Fzd4 was recomputed with the new TPM
# pipeline. The refreshed run now also reports ECs as a sending cluster
# (in addition to the original FAPs), so the table grows from 3 data rows
# to 6 data rows (header + 6 => A1:T7), and every numeric column for the
# pre-existing FAPs rows is recalculated with the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..T = Sending cluster, Ligand symbol, Receptor symbol, Target
# cluster, then the 16 NATMI metric columns, in sheet column order.
$rows = @(
    @("ECs",  "Wnt2", "Fzd4", "ECs",   1, 0.3333333333333333, 0.01070233333333333, 0.032107,  0.004227647500550067, 0.004227647500550067, 3, 1, 26.532132,         79.596396,          0.3960736634233649, 0.3960736634233648, 0.283955720708,     2.555601486372,     0.001674459833205497, 0.001674459833205497),
    @("ECs",  "Wnt2", "Fzd4", "FAPs",  1, 0.3333333333333333, 0.01070233333333333, 0.032107,  0.004227647500550067, 0.004227647500550067, 3, 1, 16.78189033333333, 50.345671,          0.2505213219764053, 0.2505213219764053, 0.1796053843107778, 1.616448458797,     0.001059115840688048, 0.001059115840688048),
    @("ECs",  "Wnt2", "Fzd4", "MuSCs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107,  0.004227647500550067, 0.004227647500550067, 3, 1, 23.67385,          71.02154999999999,  0.3534050146002298, 0.3534050146002298, 0.2533654339833333, 2.28028890585,      0.001494071826656521, 0.001494071826656521),
    @("FAPs", "Wnt2", "Fzd4", "ECs",   3, 1,                  2.520808,            7.562424,  0.99577235249945,    0.99577235249945,    3, 1, 26.532132,         79.596396,          0.3960736634233649, 0.3960736634233648, 66.88241060265601,  601.941695423904,   0.3943992035901594,   0.3943992035901593),
    @("FAPs", "Wnt2", "Fzd4", "FAPs",  3, 1,                  2.520808,            7.562424,  0.99577235249945,    0.99577235249945,    3, 1, 16.78189033333333, 50.345671,          0.2505213219764053, 0.2505213219764053, 42.30392340738933,  380.7353106665039,  0.2494622061357173,   0.2494622061357173),
    @("FAPs", "Wnt2", "Fzd4", "MuSCs", 3, 1,                  2.520808,            7.562424,  0.99577235249945,    0.99577235249945,    3, 1, 23.67385,          71.02154999999999,  0.3534050146002298, 0.3534050146002298, 59.6772304708,      537.0950742371999,  0.3519109427735733,   0.3519109427735733)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
